# Generate Report for Handback
# - Overview sheet: update handoff status text for both localization columns.
# - zh-cn / de-de sheets: add "Latest Target File" (F) and "Latest Handback
#   File" (G) columns (value + hyperlink, mirroring the existing Source/
#   Handoff-file hyperlinks) and refresh the "Latest Handback DateTime" (H).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: handoff -> handback status text. Every cell that used
# to read "Ready for handoff" (B2, C2, B3, C3 all share that string)
# now reads the handback status, so all four must be updated.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Helper: populate the Latest-Target-File (F) / Latest-Handback-File (G)
# columns for each data row of a language sheet, plus the Latest
# Handback DateTime (H) column.
# ---------------------------------------------------------------------

$mdFile = "d8671e60-fb43-46ec-b607-b1158c06422d.md"

function Apply-LanguageSheet($SheetName, $XlfFile, $MdAddress, $XlfAddress, $HandbackDateTime) {
    $ws = $wb.Worksheets.Item($SheetName)

    foreach ($row in 2, 3) {
        $fCell = $ws.Cells.Item($row, 6)   # F: Latest Target File
        $gCell = $ws.Cells.Item($row, 7)   # G: Latest Handback File

        $fCell.Value = $mdFile
        $gCell.Value = $XlfFile

        $ws.Hyperlinks.Add($fCell, $MdAddress, "", "", $mdFile)
        $ws.Hyperlinks.Add($gCell, $XlfAddress, "", "", $XlfFile)

        $fCell.Font.Underline = 2
        $fCell.Font.Color = 15570276
        $gCell.Font.Underline = 2
        $gCell.Font.Color = 15570276

        $ws.Cells.Item($row, 8).Value = $HandbackDateTime   # H: Latest Handback DateTime
    }
}

Apply-LanguageSheet "zh-cn" "d8671e60-fb43-46ec-b607-b1158c06422d.81bf13d9d024ef5e179a89dfe1e41c64b4cf54d0.zh-cn.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/6b1d5f4ef9db7586b27c66a6d7f39a29f7a8fef6/e2e/d8671e60-fb43-46ec-b607-b1158c06422d.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c07ddd45e7f31d6e02bdf0b52ed491ce8ec7f967/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d8671e60-fb43-46ec-b607-b1158c06422d.81bf13d9d024ef5e179a89dfe1e41c64b4cf54d0.zh-cn.xlf" `
    "2016-03-13 23:16:20"

Apply-LanguageSheet "de-de" "d8671e60-fb43-46ec-b607-b1158c06422d.81bf13d9d024ef5e179a89dfe1e41c64b4cf54d0.de-de.xlf" `
    "https://github.com/OpenLocalizationTest/oltest/blob/6b1d5f4ef9db7586b27c66a6d7f39a29f7a8fef6/e2e/d8671e60-fb43-46ec-b607-b1158c06422d.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8620fa6f545238f22c03290b63b3a5d3f1cdea0a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d8671e60-fb43-46ec-b607-b1158c06422d.81bf13d9d024ef5e179a89dfe1e41c64b4cf54d0.de-de.xlf" `
    "2016-03-13 23:16:26"
